$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B narrows from 14.42578125 to 13.7109375 (character width units in the
# saved OOXML). The COM ColumnWidth property is the standard way to resize a
# column; set it to the value that yields the closest column width.
$ws.Range("B:B").ColumnWidth = 12.833333333333332

# Updated data values (row 1-3, columns A & B). Row 4 is unchanged.
$ws.Range("A1").Value = -0.035865884764190649
$ws.Range("B1").Value = 0.035865884225674875

$ws.Range("A2").Value = 0.01577272595754594
$ws.Range("B2").Value = -0.01577272649659444

$ws.Range("A3").Value = 0.009363491401109373
$ws.Range("B3").Value = -0.009363491980422578
